$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update C51: change completion status from "未完成、完成一半" to "已完成"
$ws.Range("C51").Value = "已完成"

# 2. Build the new block (rows 59-66) by duplicating the existing
#    "日期：2017.9.20 四 周三 上午" block (rows 49-56), which carries the
#    same layout/styles/merge structure, then edit its text.
$ws.Range("A49:D56").Copy()
$ws.Range("A59").PasteSpecial(-4104)
$ws.Range("A49:D56").Copy()
$ws.Range("A59").PasteSpecial(-4122)

# Row heights for the new rows
$ws.Rows.Item(59).RowHeight = 22.5
$ws.Rows.Item(60).RowHeight = 22.5
$ws.Rows.Item(61).RowHeight = 22.5
$ws.Rows.Item(62).RowHeight = 22.5
$ws.Rows.Item(63).RowHeight = 22.5
$ws.Rows.Item(64).RowHeight = 22.5
$ws.Rows.Item(65).RowHeight = 22.5
$ws.Rows.Item(66).RowHeight = 22.5

# 3. Fill in the new block's content
$ws.Range("A59").Value = "日期：2017.9.20 四 周三 下午"

$ws.Range("B61").Value = "设计并修改数据库"
$ws.Range("C61").Value = "修改了一部分"

$ws.Range("B62").Value = "学习React Native框架"
$ws.Range("C62").Value = "学习中"

$ws.Range("B63").Value = "学习React Native框架"
$ws.Range("C63").Value = "学习中"

$ws.Range("B64").Value = "设计并修改数据库"
$ws.Range("C64").Value = "修改了一部分"

$ws.Range("B65").Value = "学习React Native框架"
$ws.Range("C65").Value = "学习中"

$ws.Range("A66").Value = "总结："

# 4. View state (cosmetic)
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("C64").Select()
